# Update "想去人数" (want-to-go count) values for a handful of events on the
# "展览" and "全部类型" sheets, matching the latest scraped data.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 45
$ws1.Range("F4").Value = 255
$ws1.Range("F5").Value = 3002
$ws1.Range("F6").Value = 2021
$ws1.Range("F9").Value = 1121
$ws1.Range("F11").Value = 707

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 45
$ws4.Range("F4").Value = 255
$ws4.Range("F5").Value = 3002
$ws4.Range("F6").Value = 2021
$ws4.Range("F10").Value = 1121
$ws4.Range("F12").Value = 707
